# Applies the "asset_list.xlsx" naming-convention update described in the commit:
#   "changed the names of all the assets to align with our naming convention"
#
# Net data changes on Sheet1:
#   - F,007 "Bedside Table" -> "Side Desk"
#   - New row F,009 "Wardrobe" inserted after F,008 "Vanity Stool"
#   - Seven new D-group rows (015 Bulb, 016 Blush, 017 Makeup Tubes, 018 Moomin Toy,
#     019 Pen Holder, 020 Photo Frame, 021 Poster) inserted after D,014 "Fairy Lights"
#   - The stray "Complete" note in column D next to "Character 1 body" is cleared
#   - Selection / active cell cosmetic update

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 0) Eight brand new rows are about to be inserted further up the sheet.
#    The sheet's used range stays anchored at row 300 (A1:G300) in the
#    authored workbook, so first trim eight of the trailing blank rows to
#    make room, keeping the overall row count constant.
# ---------------------------------------------------------------------------
$ws.Range("A293:A300").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 1) Rename "Bedside Table" -> "Side Desk" (row 11, column C)
# ---------------------------------------------------------------------------
$ws.Range("C11").Value2 = "Side Desk"

# ---------------------------------------------------------------------------
# 2) Insert a new row after "Vanity Stool" (row 12) for "Wardrobe"
# ---------------------------------------------------------------------------
$ws.Range("A13").EntireRow.Insert()
$ws.Range("A13").Value2 = "F"
$ws.Range("B13").Value2 = "'009"
$ws.Range("C13").Value2 = "Wardrobe"

# ---------------------------------------------------------------------------
# 3) Insert seven new rows after "Fairy Lights" (now row 27) for the new
#    Decor (D) assets, keeping the 0xx numbering scheme going.
# ---------------------------------------------------------------------------
$ws.Range("A28:A34").EntireRow.Insert()

$newDecor = @(
  @("'015", "Bulb"),
  @("'016", "Blush"),
  @("'017", "Makeup Tubes"),
  @("'018", "Moomin Toy"),
  @("'019", "Pen Holder"),
  @("'020", "Photo Frame"),
  @("'021", "Poster")
)

$r = 28
foreach ($item in $newDecor) {
    $ws.Range("A$r").Value2 = "D"
    $ws.Range("B$r").Value2 = $item[0]
    $ws.Range("C$r").Value2 = $item[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 4) Clear the leftover "Complete" marker that used to sit next to
#    "Character 1 body" (originally row 27, now shifted to row 35).
# ---------------------------------------------------------------------------
$ws.Range("D35").ClearContents()

# ---------------------------------------------------------------------------
# 5) Cosmetic: update the saved selection to match the authored file.
# ---------------------------------------------------------------------------
$ws.Range("C14").Select()
